# Update the four "(mean) X" row-label cells in the DS summary-stats table
# with their descriptive labels.
$d = $word.ActiveDocument

$replacements = @(
    @{old = "(mean) war";      new = "War involving own nation"},
    @{old = "(mean) neutral";  new = "Neutrality of own nation"},
    @{old = "(mean) TONMOD";   new = "Tonnage standardized on British measured tons, 1773-1835"},
    @{old = "(mean) crowd";    new = "Number of embarked enslaved persons per ton"}
)

foreach ($r in $replacements) {
    $found = $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $r.new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $($r.old)"
    }
}
